$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4321.310770455058,
    4321.310770455058,
    4302.643575320226,
    4176.633223453386,
    4176.633223453386,
    4176.633223453386,
    4176.633223453386,
    4176.633223453386,
    4176.633223453386,
    4176.633223453386,
    3875.242430188273
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
